$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BD = New-Object 'double[,]' 24,3
$BD[0,0] = 0.5454946546819031
$BD[0,1] = 0.2882628619767331
$BD[0,2] = 0.2026351672010023
$BD[1,0] = 0.5097306866138069
$BD[1,1] = 0.2892496447514041
$BD[1,2] = 0.2010410331994308
$BD[2,0] = 0.4878859361623711
$BD[2,1] = 0.2899043686539642
$BD[2,2] = 0.2001307733999838
$BD[3,0] = 0.4790135267736844
$BD[3,1] = 0.2901834899364921
$BD[3,2] = 0.1997771561068475
$BD[4,0] = 0.4775420730959183
$BD[4,1] = 0.2902305827276273
$BD[4,2] = 0.1997194868615253
$BD[5,0] = 0.4877661593318123
$BD[5,1] = 0.2899080830663152
$BD[5,2] = 0.2001259341478629
$BD[6,0] = 0.5331398996016787
$BD[6,1] = 0.2885929904168414
$BD[6,2] = 0.2020713332902346
$BD[7,0] = 0.6229989431832905
$BD[7,1] = 0.2864000069668577
$BD[7,2] = 0.206426747585823
$BD[8,0] = 0.6895261261655889
$BD[8,1] = 0.2850219652229171
$BD[8,2] = 0.20995213156462
$BD[9,0] = 0.7198954785399962
$BD[9,1] = 0.2844452401243984
$BD[9,2] = 0.2116257660260743
$BD[10,0] = 0.7314101361463941
$BD[10,1] = 0.2842340261781757
$BD[10,2] = 0.2122695028664339
$BD[11,0] = 0.7289296185216187
$BD[11,1] = 0.2842791960936921
$BD[11,2] = 0.2121304206776671
$BD[12,0] = 0.7208425113073815
$BD[12,1] = 0.2844277197280789
$BD[12,2] = 0.211678527368619
$BD[13,0] = 0.7158907826549239
$BD[13,1] = 0.284519628740135
$BD[13,2] = 0.2114030252022019
$BD[14,0] = 0.6875434806246687
$BD[14,1] = 0.285060661976992
$BD[14,2] = 0.2098441551763131
$BD[15,0] = 0.6701798784419282
$BD[15,1] = 0.285405390418866
$BD[15,2] = 0.2089056853891549
$BD[16,0] = 0.6602027984388883
$BD[16,1] = 0.2856083916688377
$BD[16,2] = 0.2083724887257858
$BD[17,0] = 0.6568264724272126
$BD[17,1] = 0.2856779365085274
$BD[17,2] = 0.2081930913451657
$BD[18,0] = 0.6720272336774826
$BD[18,1] = 0.2853682049588784
$BD[18,2] = 0.2090049060665393
$BD[19,0] = 0.7232175048944214
$BD[19,1] = 0.2843839001509707
$BD[19,2] = 0.2118109895547917
$BD[20,0] = 0.756757133903136
$BD[20,1] = 0.2837824318054771
$BD[20,2] = 0.2137029871204845
$BD[21,0] = 0.7388489870056958
$BD[21,1] = 0.2840996298593836
$BD[21,2] = 0.2126879085242024
$BD[22,0] = 0.6711920272292957
$BD[22,1] = 0.2853850015134753
$BD[22,2] = 0.2089600286369944
$BD[23,0] = 0.5985984001499105
$BD[23,1] = 0.2869521730433462
$BD[23,2] = 0.2051910468296967
$ws.Range("B2:D25").Value = $BD

$FH = New-Object 'double[,]' 24,3
$FH[0,0] = 1.652177775467138
$FH[0,1] = 0.9646673426281538
$FH[0,2] = 1.016199478200917
$FH[1,0] = 1.661741970935012
$FH[1,1] = 0.9717949163800839
$FH[1,2] = 1.023297572375746
$FH[2,0] = 1.668316867349802
$FH[2,1] = 0.9766533506422093
$FH[2,2] = 1.028006955142949
$FH[3,0] = 1.671173039138402
$FH[3,1] = 0.9787544865351805
$FH[3,2] = 1.03001447093061
$FH[4,0] = 1.671657991749392
$FH[4,1] = 0.9791107056769732
$FH[4,2] = 1.030353160051995
$FH[5,0] = 1.66835467034646
$FH[5,1] = 0.9766811961360631
$FH[5,2] = 1.028033671111196
$FH[6,0] = 1.655329845529948
$FH[6,1] = 0.9670249127210582
$FH[6,2] = 1.018574092006567
$FH[7,0] = 1.635352854350998
$FH[7,1] = 0.9519115287446454
$FH[7,2] = 1.002804987931206
$FH[8,0] = 1.624056833362147
$FH[8,1] = 0.9431351053535764
$FH[8,2] = 0.99290820832978
$FH[9,0] = 1.619649814029437
$FH[9,1] = 0.9396472590771907
$FH[9,2] = 0.9887712115113487
$FH[10,0] = 1.618085994503573
$FH[10,1] = 0.9383990111099649
$FH[10,2] = 0.9872570313326605
$FH[11,0] = 1.618418122955433
$FH[11,1] = 0.9386646187120817
$FH[11,2] = 0.9875808075717316
$FH[12,0] = 1.619519053867862
$FH[12,1] = 0.9395431116907957
$FH[12,2] = 0.9886455891212051
$FH[13,0] = 1.620207077603496
$FH[13,1] = 0.9400906576561709
$FH[13,2] = 0.9893046213289551
$FH[14,0] = 1.624359548686989
$FH[14,1] = 0.9433731949310555
$FH[14,2] = 0.9931859098089006
$FH[15,0] = 1.627094206982036
$FH[15,1] = 0.9455161379956465
$FH[15,2] = 0.9956604007987266
$FH[16,0] = 1.628735982261027
$FH[16,1] = 0.9467961990118923
$FH[16,2] = 0.9971180326523736
$FH[17,0] = 1.629303693494663
$FH[17,1] = 0.9472377643993681
$FH[17,2] = 0.9976174675948855
$FH[18,0] = 1.626795971229853
$FH[18,1] = 0.9452831025630744
$FH[18,2] = 0.9953934304857768
$FH[19,0] = 1.619192834910123
$FH[19,1] = 0.9392831090252116
$FH[19,2] = 0.9883314152715599
$FH[20,0] = 1.614835816949252
$FH[20,1] = 0.9357844730639897
$FH[20,2] = 0.9840214285232349
$FH[21,0] = 1.617105294327665
$FH[21,1] = 0.937613095282714
$FH[21,2] = 0.9862938306049216
$FH[22,0] = 1.626930586852495
$FH[22,1] = 0.9453883081935501
$FH[22,2] = 0.9955140186911677
$FH[23,0] = 1.64016259217842
$FH[23,1] = 0.9555911726262494
$FH[23,2] = 1.006773901913739
$ws.Range("F2:H25").Value = $FH

$JM = New-Object 'double[,]' 24,4
$JM[0,0] = 0.2691701291591979
$JM[0,1] = 0.2568629249896333
$JM[0,2] = 0.3037962595540051
$JM[0,3] = 0.1893672929400942
$JM[1,0] = 0.2710674796308936
$JM[1,1] = 0.225492109653743
$JM[1,2] = 0.3013966425616843
$JM[1,3] = 0.1826640379586415
$JM[2,0] = 0.2723060814417337
$JM[2,1] = 0.206187700513496
$JM[2,2] = 0.3000238292334458
$JM[2,3] = 0.1786061342751992
$JM[3,0] = 0.2728293615336348
$JM[3,1] = 0.198310821415518
$JM[3,2] = 0.2994898011521201
$JM[3,3] = 0.1769672310893959
$JM[4,0] = 0.2729173723524259
$JM[4,1] = 0.1970022728861096
$JM[4,2] = 0.2994026640237095
$JM[4,3] = 0.1766959865760143
$JM[5,0] = 0.2723130634768882
$JM[5,1] = 0.2060815106553093
$JM[5,2] = 0.3000165241256028
$JM[5,3] = 0.1785839716173925
$JM[6,0] = 0.2698090720383792
$JM[6,1] = 0.246055422086414
$JM[6,2] = 0.3029480693546205
$JM[6,3] = 0.1870440828016164
$JM[7,0] = 0.2654816495926422
$JM[7,1] = 0.3240862564109648
$JM[7,2] = 0.3094900807662029
$JM[7,3] = 0.2040878089891258
$JM[8,0] = 0.2626558749904859
$JM[8,1] = 0.3811761331309071
$JM[8,2] = 0.3147746027571543
$JM[8,3] = 0.2168794590789673
$JM[9,0] = 0.2614467583220943
$JM[9,1] = 0.4070916649524463
$JM[9,2] = 0.3172813695602912
$JM[9,3] = 0.2227558966981533
$JM[10,0] = 0.2609998476911279
$JM[10,1] = 0.4168968303930285
$JM[10,2] = 0.3182452922550283
$JM[10,3] = 0.2249892676798382
$JM[11,0] = 0.2610956109918696
$JM[11,1] = 0.4147854991877011
$JM[11,2] = 0.3180370439519606
$JM[11,3] = 0.2245079141064394
$JM[12,0] = 0.261409771275547
$JM[12,1] = 0.4078985151857069
$JM[12,2] = 0.3173603788463595
$JM[12,3] = 0.2229394762779506
$JM[13,0] = 0.2616036296051423
$JM[13,1] = 0.4036789164550498
$JM[13,2] = 0.3169478081348984
$JM[13,3] = 0.2219798115449549
$JM[14,0] = 0.2627364278997426
$JM[14,1] = 0.3794813367191807
$JM[14,2] = 0.3146128382485216
$JM[14,3] = 0.2164965616087002
$JM[15,0] = 0.2634509001346981
$JM[15,1] = 0.3646224181659647
$JM[15,2] = 0.3132066571302374
$JM[15,3] = 0.2131473592034112
$JM[16,0] = 0.2638690326979738
$JM[16,1] = 0.3560708248620585
$JM[16,2] = 0.3124075431689874
$JM[16,3] = 0.2112264043453607
$JM[17,0] = 0.2640118402612437
$JM[17,1] = 0.3531745379441134
$JM[17,2] = 0.3121386438337765
$JM[17,3] = 0.210576937727069
$JM[18,0] = 0.2633740996774865
$JM[18,1] = 0.3662047112391349
$JM[18,2] = 0.3133553460971541
$JM[18,3] = 0.2135033280285512
$JM[19,0] = 0.2613171976303157
$JM[19,1] = 0.4099216242296109
$JM[19,2] = 0.3175587347525521
$JM[19,3] = 0.2233999462373504
$JM[20,0] = 0.2600367406710866
$JM[20,1] = 0.4384434327694748
$JM[20,2] = 0.3203913127858584
$JM[20,3] = 0.229915043619414
$JM[21,0] = 0.260714309939523
$JM[21,1] = 0.4232255488001613
$JM[21,2] = 0.3188717360886955
$JM[21,3] = 0.2264335613130015
$JM[22,0] = 0.2634087981961279
$JM[22,1] = 0.3654893845406093
$JM[22,2] = 0.3132880947814414
$JM[22,3] = 0.2133423803492391
$JM[23,0] = 0.2665901121760132
$JM[23,1] = 0.3240862564109648
$JM[23,2] = 0.3076359370707777
$JM[23,3] = 0.1994291885112993
$ws.Range("J2:M25").Value = $JM

$OCol = New-Object 'double[,]' 24,1
$OCol[0,0] = 4.008469675755563
$OCol[1,0] = 4.038187385106966
$OCol[2,0] = 4.058182405032895
$OCol[3,0] = 4.066770503190369
$OCol[4,0] = 4.068223133729504
$OCol[5,0] = 4.058296445248459
$OCol[6,0] = 4.018353712961527
$OCol[7,0] = 3.953882235886425
$OCol[8,0] = 3.914942282056245
$OCol[9,0] = 3.899053281728101
$OCol[10,0] = 3.893298632722889
$OCol[11,0] = 3.894526343374793
$OCol[12,0] = 3.898574590488067
$OCol[13,0] = 3.901088393385777
$OCol[14,0] = 3.916017366694234
$OCol[15,0] = 3.925643044955791
$OCol[16,0] = 3.931351271762736
$OCol[17,0] = 3.93331349249614
$OCol[18,0] = 3.924600597198975
$OCol[19,0] = 3.897378409013214
$OCol[20,0] = 3.881115145715427
$OCol[21,0] = 3.88965542808333
$OCol[22,0] = 3.92507134489091
$OCol[23,0] = 3.969842136919127
$ws.Range("O2:O25").Value = $OCol
